# Update the dSF (column F) values on the active worksheet to reflect
# a repull of data / recalculated mean values, per the commit message
# "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -6
    4  = -1
    8  = -8
    9  = -3
    12 = -2
    14 = -6
    15 = -6
    16 = 4
    17 = -4
    18 = -2
    19 = -4
    20 = 2
    21 = 2
    22 = 3
    23 = -4
    24 = 1
    25 = -4
    26 = -2
    27 = -4
    28 = -2
    29 = -3
    30 = -3
    31 = -4
    32 = 2
    34 = -4
    35 = -1
    37 = -1
    38 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
